# Generate Report for Handback
#
# This script fills in the "handback" columns (Latest Target File,
# Latest Handback File, Latest Handback DateTime) for both the zh-cn and
# de-de language sheets, updates the "Status" text everywhere it appears
# (Overview sheet + per-language sheets) from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens a few columns so the new,
# longer text/hyperlinks are not clipped.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Hyperlink targets re-used from the existing "Source File Name" links.
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f18abe53d2f88cd3ed20620061cc82a66666e7c/e2e/3a73c11f-eca7-41bf-9da7-aa9e86668101.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f18abe53d2f88cd3ed20620061cc82a66666e7c/e2e/5c98e9ba-e5ba-4b2b-b50c-fbfac42e90b1.md"
$mdName1 = "3a73c11f-eca7-41bf-9da7-aa9e86668101.md"
$mdName2 = "5c98e9ba-e5ba-4b2b-b50c-fbfac42e90b1.md"

# Color/format used by the workbook's existing custom "HyperLink" cell style
# (underline, font color #6495ED - cornflowerblue). OLE colors are BGR, so
# RGB(0x64,0x95,0xED) becomes 0xED9564.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: refresh the Status column shown for each language
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$languages = @(
    @{
        Sheet = "zh-cn"
        Xlf1 = "3a73c11f-eca7-41bf-9da7-aa9e86668101.f437634fb4767deb3fcebbc99ce22a9882f0cda6.zh-cn.xlf"
        Xlf2 = "5c98e9ba-e5ba-4b2b-b50c-fbfac42e90b1.2a87cd170647c6eed6c6da96a32c062fb419b44e.zh-cn.xlf"
        HandbackDateTime = "2016-08-31 13:08:08"
    },
    @{
        Sheet = "de-de"
        Xlf1 = "3a73c11f-eca7-41bf-9da7-aa9e86668101.f437634fb4767deb3fcebbc99ce22a9882f0cda6.de-de.xlf"
        Xlf2 = "5c98e9ba-e5ba-4b2b-b50c-fbfac42e90b1.2a87cd170647c6eed6c6da96a32c062fb419b44e.de-de.xlf"
        HandbackDateTime = "2016-08-31 13:08:30"
    }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column for both rows.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 (3a73c11f... file): Latest Target File / Latest Handback File /
    # Latest Handback DateTime.
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, "", "", $mdName1)
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("J2").Value = $lang.Xlf1
    $ws.Range("K2").Value = $lang.HandbackDateTime

    # Row 3 (5c98e9ba... file).
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, "", "", $mdName2)
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = $hyperlinkColor
    $ws.Range("J3").Value = $lang.Xlf2
    $ws.Range("K3").Value = $lang.HandbackDateTime

    # Widen columns so the new content is fully visible.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Write-Output "Handback report generated."
